# Weekly refresh: shift existing rows 151-253 down into 152-254 (a new
# week's worth of "Segunda" observation is inserted at row 151, pushing
# the historical rows down by one and dropping the oldest row that was
# at 254).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19)   # D, L, M, N, O, P, Q, R, S

for ($r = 254; $r -ge 152; $r--) {
    foreach ($c in $cols) {
        $src = $ws.Cells.Item($r - 1, $c).Value2
        $ws.Cells.Item($r, $c).Value = $src
    }
}

# New data for row 151 (the newly reported week).
$ws.Cells.Item(151, 4).Value  = 44810               # Fecha
$ws.Cells.Item(151, 12).Value = "Segunda"           # Calidad
$ws.Cells.Item(151, 13).Value = 120                 # Volumen
$ws.Cells.Item(151, 14).Value = 8000                # Precio minimo
$ws.Cells.Item(151, 15).Value = 9000                # Precio maximo
$ws.Cells.Item(151, 16).Value = 8500                # Precio promedio ponderado
$ws.Cells.Item(151, 17).Value = '$/bandeja 7 kilos'  # Unidad de comercializacion
$ws.Cells.Item(151, 18).Value = 'Provincia de Melipilla' # Origen
$ws.Cells.Item(151, 19).Value = 1214                # Precio $/Kg
